$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2:D358").Value = 2
$ws.Range("E7").Select()
